$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("L22").Value = 1612.02
$ws1.Range("M22").Value = 10113.13
$ws1.Range("L26").Value = "1 de 24"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F22").Value = 11725.15
$ws2.Range("F26").Value = 21202.47

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D11").Value = 219.99
$ws3.Range("E11").Value = 2702.23458185274
$ws3.Range("F11").Value = 0.07528168826111326

$ws3.Range("D12").Value = 20982.48
$ws3.Range("E12").Value = 6972.5
$ws3.Range("F12").Value = 0.7505811129179846

$ws3.Range("D14").Value = 21202.47
$ws3.Range("E14").Value = 21000.91110009469
$ws3.Range("F14").Value = 0.5023879473000904
